$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resoconto finale")
$ws.Range("T4").Value = 7
$r = $ws.Range("T4")
Write-Output $r.Value2
